# TC05_Canine_Filter_Diagnosis-OsteoSarcoma.xlsx — "startup" sheet fix-ups.
#
# Commit message: "Fixed Diagnosis, FileAssociation, FileFormat, FileType,
# NeuteredStatus, PrimeDiseaseSite"
#
# The real content change is in the Cypher query stored in B2 (the "Cases"
# tab query, row 2 / CasesTab): the final RETURN clause no longer projects
# coalesce(co.cohort_description,'') AS "Cohort", so that trailing line
# (and the now-dangling trailing comma on the previous line) is removed.
# Everything else the query does (including the still-present
# "MATCH (co:cohort)..." clause) is left untouched, matching the diff
# exactly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$casesQuery = @'
MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
MATCH (c)<--(diag:diagnosis)
MATCH (co:cohort)<-[*]-(c)
WHERE diag.disease_term IN ['Osteosarcoma']
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS `Case ID` ,
        coalesce(s.clinical_study_designation, '') AS `Study Code` ,
        coalesce(s.clinical_study_type, '') AS  `Study Type`,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term,'') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS `Stage of Disease` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS `Neutered Status`,
        coalesce(demo.weight, '') AS `Weight (kg)`,
        coalesce(diag.best_response, '') AS `Response to Treatment`
'@

$ws.Range("B2").Value = $casesQuery

# The workbook was also re-saved with a different selection/view: the
# previously-selected cell B4 (scrolled so row 4 was the top-left row) is
# now B2, with the window scrolled back to the top of the sheet.
$null = $ws.Range("B2").Select()

# ... and the zoom was changed to 100%.
$excel.ActiveWindow.Zoom = 100
